$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header / row-label text: insert line breaks so long labels wrap onto
#    multiple lines (cells already have wrapText formatting).
# ---------------------------------------------------------------------------
$ws.Range("B1").Value = "Стоимость`nТС"
$ws.Range("C1").Value = "Стоимость`nобслуживания`nТС"
$ws.Range("I1").Value = "Год`nвыпуска"
$ws.Range("J1").Value = "Вектор`nприоритетов"
$ws.Range("A2").Value = "Стоимость`nТС"
$ws.Range("A3").Value = "Стоимость`nобслуживания`nТС"
$ws.Range("A9").Value = "Год`nвыпуска"

# Re-run AutoFit row by row so the row heights stay on the default (the
# workbook does not carry any explicit/custom row heights).
$ws.Rows.Item(1).AutoFit()
$ws.Rows.Item(2).AutoFit()
$ws.Rows.Item(3).AutoFit()
$ws.Rows.Item(9).AutoFit()

# ---------------------------------------------------------------------------
# 2. Updated priority-vector / consistency numbers in column J.
# ---------------------------------------------------------------------------
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "0.328"
$ws.Range("J3").NumberFormat = "@"
$ws.Range("J3").Value = "0.159"
$ws.Range("J4").NumberFormat = "@"
$ws.Range("J4").Value = "0.232"
$ws.Range("J5").NumberFormat = "@"
$ws.Range("J5").Value = "0.107"
$ws.Range("J6").NumberFormat = "@"
$ws.Range("J6").Value = "0.048"
$ws.Range("J7").NumberFormat = "@"
$ws.Range("J7").Value = "0.071"
$ws.Range("J8").NumberFormat = "@"
$ws.Range("J8").Value = "0.033"
$ws.Range("J9").NumberFormat = "@"
$ws.Range("J9").Value = "0.023"

# The NumberFormat tweak above was only needed to stop Excel from turning
# the numeric-looking text into a real number; restore the original
# (General / bordered / wrap-text) cell formatting by copying it back from
# the untouched neighbour column.
$ws.Range("I2:I9").Copy()
$ws.Range("J2:J9").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("J10").Value = "λ_max = 8.288"
$ws.Range("J11").Value = "ИС = 0.041"
$ws.Range("J12").Value = "ОС = 0.029"

# ---------------------------------------------------------------------------
# 3. Column widths.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 16.0
$ws.Columns.Item(2).ColumnWidth = 12.333333333333334
$ws.Columns.Item(3).ColumnWidth = 16.0
$ws.Columns.Item(9).ColumnWidth = 10.0
$ws.Columns.Item(10).ColumnWidth = 14.833333333333334
